$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Valor Mora" total and the period count -------------------
$ws.Range("E11").Value = 660000
$ws.Range("F13").Value = 5

# --- Add a new worker/period row (2507) and renumber the existing ones ----
# Before: rows 16-19 hold periods 2503,2504,2505,2506 (oldest->newest, row19
# styled as the last/bottom row of the mini-table).
# After:  rows 16-20 hold periods 2507,2506,2505,2504,2503 (newest->oldest),
# row20 becomes the new last/bottom row.

# 1) Insert a blank row at 20 (pushes the signature block down from 24/25 to 25/26)
$ws.Rows("20").Insert()

# 2) Give the new row20 the "bottom of table" formatting that row19 used to have
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# 3) Row19 is no longer the last row of the table, so restyle it like rows 16-18
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# 4) Fill in the new row20 data (copy of the worker row, with period 2503)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143392800"
$ws.Range("D20").Value = "JORGE ENRIQUE CIRO TORO"
$ws.Range("E20").Value = "2503"
$ws.Range("F20").Value = 132000
$ws.Range("G20").Value = 3300000

# 5) Shift the period values: row16=2507, row17=2506, row18=2505, row19=2504
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2504"

"done"
